$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(68, "...nights", "～泊|～はく"),
    @(69, "with...", "～付|～つき"),
    @(70, "with breakfast", "朝食付き|ちょうしょくつき"),
    @(71, "one night with two meals", "一泊二食付|いっぱくにしょくつき"),
    @(72, "checking in", "チェックイン（する）"),
    @(73, "checking out", "チェックアウト（する）"),
    @(74, "single room", "シングル"),
    @(75, "double room", "ダブル"),
    @(76, "twin room", "ツイン"),
    @(77, "...person(s)", "～名|～めい"),
    @(78, "receptionist; front desk", "フロント"),
    @(79, "non-smoking", "禁煙|きんえん"),
    @(80, "smoking", "喫煙|きつえん"),
    @(81, "large spa", "大浴場|だいよくじょう"),
    @(82, "I have a reservation under the name of XX.", "予約した○○です。|よやくしたまるまるです。"),
    @(83, "Could you check me in, please?", "チェックインお願いします。|チェックインおねがいします。"),
    @(84, "Can I pay by XX?", "○○で払えますか。|まるまるではらえますか。"),
    @(85, "Can you change US dollars to Japanese yen?", "アメリカドルを日本円に両替できますか。|アメリカドルをにほんえんにりょうがえできますか。"),
    @(86, "Would you call a taxi for us?", "タクシーを呼んでもらえますか。|タクシーをよんでもらえますか。"),
    @(87, "Are there any restaurants you recommend near here?", "近くにおすすめのレストランがありますか。|ちかくにおすすめのレストランがありますか。"),
    @(88, "Could you keep my luggage until 2 o'clock?", "二時まで荷物を預かってくれませんか。|にじまでにもつをあずかってくれませんか。"),
)

foreach ($row in $newRows) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}
